# Camden Scholl Programming List
# Populates Sheet1 with a title + a "Title/Description/Skills/github link"
# table, formatted with the built-in Heading 1 / Heading 2 / Explanatory Text
# / Hyperlink cell styles, wrapped text, column widths, row heights,
# hyperlinks, page orientation, and the current selection -- reproducing the
# authored workbook described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Title (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Camden Scholl Programming List"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 14

# ---------------------------------------------------------------------
# Header row (row 4)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Title"
$ws.Range("A4").Style = "Heading 1"

$ws.Range("B4").Value = "Description"
$ws.Range("B4").Style = "Heading 2"

$ws.Range("C4").Value = "Skills"
$ws.Range("C4").Style = "Heading 2"

$ws.Range("D4").Value = "github link"
$ws.Range("D4").Style = "Heading 2"

# ---------------------------------------------------------------------
# Data rows (5-8)
# ---------------------------------------------------------------------
$titles = @("heartBeat", "CABERDLE", "YouTube Homepage replica", "Discord Bot")

$descriptions = @(
    "Full stack project using Vue.js and Firebase to output Spotify-listed songs based on a user’s heartbeat and mood.",
    "Worked with partner to create a more player-friendly version of WORDLE with added difficulty levels and the ability to replay.",
    "Practicing a variety of basic HTML and CSS skills. Followed a tutorial by SuperSimpleDev (link in README)",
    "Created a Discord bot that responds to various commands. Followed the instructions on discord.py"
)

$skills = @(
    "JavaScript, HTML, CSS, node.js, Vue, Firebase ",
    "Java",
    "HTML, CSS",
    "JavaScript, node.js"
)

$links = @(
    "https://github.com/c-l-scholl/stp-vue-fb",
    "https://github.com/mac-comp128-s22/128-project-ben-and-camden2",
    "https://github.com/c-l-scholl/YT-home-page-UI-replica",
    "https://github.com/c-l-scholl/discord-bot"
)

for ($i = 0; $i -lt 4; $i++) {
    $row = 5 + $i

    $ws.Cells.Item($row, 1).Value = $titles[$i]

    $descCell = $ws.Cells.Item($row, 2)
    $descCell.Value = $descriptions[$i]
    $descCell.Style = "Explanatory Text"
    $descCell.WrapText = $true

    $ws.Cells.Item($row, 3).Value = $skills[$i]

    $linkCell = $ws.Cells.Item($row, 4)
    $linkCell.Value = $links[$i]
    $ws.Hyperlinks.Add($linkCell, $links[$i]) | Out-Null
}

# Skills cell in row 5 wraps text (others keep the default style).
$ws.Range("C5").WrapText = $true

# ---------------------------------------------------------------------
# Column widths (approximate best-fit autosize) & row heights
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.721354166666668
$ws.Columns.Item(2).ColumnWidth = 47.608072916666664
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 58.385416666666664

$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(4).RowHeight = 20.4
$ws.Rows.Item(5).RowHeight = 43.8
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 28.8

# ---------------------------------------------------------------------
# Page setup, selection
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

$ws.Range("C5").Select() | Out-Null

Write-Host "Applied Camden Scholl Programming List edits"
